$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. Expand the Bensound.com hyperlink display text to the full URL.
# ---------------------------------------------------------------------------
$d.Content.Find.Execute("Bensound.com", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "https://www.bensound.com/royalty-free-music", 2)

# ---------------------------------------------------------------------------
# 2. Move the "_GoBack" bookmark from the empty paragraph (just after
#    "animation clock") down onto the final YouTube link paragraph.
# ---------------------------------------------------------------------------
$d.Bookmarks("_GoBack").Delete()

# ---------------------------------------------------------------------------
# 3. Re-order the grammar-check markers around "wrong" / " answer  sound" so
#    the proofErr gramStart sits before "wrong" and gramEnd sits right after
#    it (instead of bracketing "answer  sound").
#    (Neither of the edits above add/remove paragraphs, so lookup by stable
#    text match -- trimming the trailing paragraph-mark "\r" that Range.Text
#    always reports -- finds the right paragraph regardless of index drift.)
# ---------------------------------------------------------------------------
$wrongPara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.TrimEnd([char]13) -eq "wrong answer  sound") {
        $wrongPara = $p
        break
    }
}

# Drop the leading run ("wrong ") entirely -- this collapses its empty <w:r>
# away, which slides the gramStart marker to the very front of the paragraph.
$r1 = $d.Range($wrongPara.Range.Start, $wrongPara.Range.Start + 6)
$r1.Text = ""

# Type "wrong" back in at the paragraph start; it merges into the remaining
# run (still positioned after gramStart).
$insA = $d.Range($wrongPara.Range.Start, $wrongPara.Range.Start)
$insA.Text = "wrong"

# Strip "answer  sound" back off that same run, leaving gramStart, "wrong",
# gramEnd with nothing in between the marker pair and the run.
$r3 = $d.Range($wrongPara.Range.Start + 5, $wrongPara.Range.End - 1)
$r3.Text = ""

# Finally append " answer  sound" via InsertAfter on the paragraph's own
# Range -- this lands the new run after gramEnd rather than merging into the
# "wrong" run that precedes it.
$wrongPara.Range.InsertAfter(" answer  sound")

# ---------------------------------------------------------------------------
# 4. Re-create the "_GoBack" bookmark around the final YouTube URL run.
# ---------------------------------------------------------------------------
$ytPara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.TrimEnd([char]13) -eq "https://www.youtube.com/watch?v=QZYFVgBu7cE") {
        $ytPara = $p
        break
    }
}
$ytTarget = $d.Range($ytPara.Range.Start, $ytPara.Range.End - 1)
$d.Bookmarks.Add("_GoBack", $ytTarget)
